{"js": "// The template's intro sentence was split across five separate runs\n// (\"A simple \" / \"demonstration\" / \" of a \" / \"query\" / \" :\"). Re-write it\n// as a single run holding the full sentence.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst introParagraph = paragraphs.items[0];\nintroParagraph\n    .getRange()\n    .insertText(\"A simple demonstration of a query :\", Word.InsertLocation.replace);\n\n// The \"anydsl\" paragraph had an explicit orange accent color on its run;\n// reset that run's font color back to automatic (no explicit color).\nconst anydslResults = body.search(\"anydsl\", { matchCase: true });\nanydslResults.load(\"text\");\nawait context.sync();\n\nconst anydslRun = anydslResults.items[0];\nconst resetRunOoxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r><w:rPr><w:color w:val=\"auto\"/></w:rPr><w:t>anydsl</w:t></w:r></w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\nanydslRun.insertOoxml(resetRunOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The template's intro sentence was split across five separate runs\n# (\"A simple \" / \"demonstration\" / \" of a \" / \"query\" / \" :\"). Locate that\n# paragraph and rewrite its content as a single run holding the full\n# sentence (same text, just no longer fragmented into multiple runs).\n$introFind = $d.Content.Find\n$introFind.Text = \"A simple\"\n$introFind.Execute() | Out-Null\n\n$introParagraph = $introFind.Paragraphs(1).Range\n$introParagraph.SetRange($introParagraph.Start, $introParagraph.End - 1)\n$introParagraph.Text = \"A simple demonstration of a query :\"\n\n# The \"anydsl\" paragraph had an explicit orange accent color on its run;\n# reset that run's font color back to automatic (no explicit color).\n$anydslFind = $d.Content.Find\n$anydslFind.Text = \"anydsl\"\n$anydslFind.Execute() | Out-Null\n$anydslFind.Parent.Font.Color = \"wdColorAutomatic\"\n"}
